$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "WARNING: replace failed for: $old"
    }
}

# 1. "...call on their representatives on the HESTA Board to divest..."
#    -> "...call on HESTA to divest..."
Replace-Text "call on their representatives on the HESTA Board to divest" "call on HESTA to divest"

# 2. "...a motion to go before [Union] National Council . " (with gramStart/gramEnd proofErr)
#    -> "...a motion to go before the [Union] National Council. "
Replace-Text "to go before [Union] National Council . " "to go before the [Union] National Council. "

# 3. "...detention facilities. We note that doctors..."
#    -> "...detention facilities, and the lack of access to medical care. We note that doctors..."
Replace-Text "the Nauru and Manus Island detention facilities. We note that doctors" "the Nauru and Manus Island detention facilities, and the lack of access to medical care. We note that doctors"

# 4. "Mandatory detention only occurs because the government supports it and because commercial investors support it. Stopping"
#    -> "Mandatory detention only occurs because of government support and commercial investment. Stopping"
Replace-Text "Mandatory detention only occurs because the government supports it and because commercial investors support it. Stopping" "Mandatory detention only occurs because of government support and commercial investment. Stopping"

# 5. " [Union's] role on the HESTA board gives members to power to have real influence ... HESTA's decision whether or not to purchase shares"
#    -> " [Union's] role on the HESTA Board gives [Union] members the power to have real influence ... HESTA's decision on whether or not to purchase shares"
Replace-Text " [Union’s] role on the HESTA board gives members to power to have real influence in the corporate decision-making of this industry super fund. HESTA’s decision whether or not to purchase shares" " [Union’s] role on the HESTA Board gives [Union] members the power to have real influence in the corporate decision-making of this industry super fund. HESTA’s decision on whether or not to purchase shares"

# 6. Delete the empty paragraph right before "[Union] officers will make clear representations..."
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "") {
        $next = $p.Next()
        if ($next -ne $null -and $next.Range.Text.StartsWith("[Union] officers will make clear representations")) {
            $p.Range.Delete()
            break
        }
    }
}

# 7. "[Union] officers will make clear representations to HESTA, calling on its board of Directors to change the fund's investment practices in keeping with the values of this statement. That implies an important role..."
#    -> "...calling on its Board of Directors... This implies an important role..."
Replace-Text "[Union] officers will make clear representations to HESTA, calling on its board of Directors to change the fund’s investment practices in keeping with the values of this statement. That implies an important role" "[Union] officers will make clear representations to HESTA, calling on its Board of Directors to change the fund’s investment practices in keeping with the values of this statement. This implies an important role"

# 8. "to call for the Board of HESTA Board to commit to a change in HESTA's practices"
#    -> "to call for a change in HESTA's practices"
Replace-Text "to call for the Board of HESTA Board to commit to a change in HESTA’s practices" "to call for a change in HESTA’s practices"

Write-Output "done"
